# 300803-指南针.xlsx: add 2022-Q4 data
#
# 1. Insert a new "2022-Q4" worksheet right after "总计" (i.e. before the
#    sheet currently named "2022-Q3"), filled with the new quarter's fund
#    holdings table.
# 2. In the "总计" summary sheet, insert a new row right under the header
#    with the 2022-Q4 totals; the existing quarters shift down automatically.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: "总计" sheet - insert the 2022-Q4 summary row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()
# Row insert drags the header's bold/border formatting onto the new row;
# strip that so the new data row starts out unstyled like its siblings.
$summary.Range("A2:D2").ClearFormats()

# Re-apply the boxed/bold/centered style used by the rest of column A to
# the new A2 cell (copy format only, not value) by pulling it from A3,
# a surviving data cell that still carries the original style.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 15
$summary.Range("D2").Value = 4.7

# ---------------------------------------------------------------------------
# Step 2: brand-new "2022-Q4" worksheet with the per-fund breakdown.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Add($q3)
$newSheet.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")

$rows = @(
  @("008188", "前海开源稳健增长三年持有期混合", "21.91", "91.82", "7.62", "1.6695", 5),
  @("004702", "南方金融主题灵活配置混合A", "12.97", "92.71", "4.89", "0.6342", 6),
  @("001178", "前海开源再融资主题精选股票", "5.94", "92.55", "7.61", "0.4520", 5),
  @("007346", "易方达科技创新混合", "17.71", "89.19", "2.55", "0.4516", 5),
  @("001103", "前海开源工业革命4.0灵活配置混合", "4.12", "87.14", "7.80", "0.3214", 5),
  @("013500", "南方金融主题灵活配置混合C", "4.80", "92.71", "4.89", "0.2347", 6),
  @("013610", "中信保诚前瞻优势混合", "12.64", "82.82", "1.81", "0.2288", 10),
  @("011287", "前海开源聚慧三年持有期混合", "2.88", "92.27", "7.61", "0.2192", 5),
  @("006775", "前海开源优质成长混合", "2.55", "92.63", "7.52", "0.1918", 5),
  @("006216", "前海开源价值成长灵活配置混合A", "1.14", "91.81", "7.51", "0.0856", 5),
  @("002407", "前海开源恒远灵活配置混合", "1.00", "92.86", "7.98", "0.0798", 5),
  @("159851", "华宝中证金融科技主题ETF", "1.98", "98.27", "3.39", "0.0671", 6),
  @("006217", "前海开源价值成长灵活配置混合C", "0.47", "91.81", "7.51", "0.0353", 5),
  @("516100", "华夏中证金融科技主题ETF", "0.60", "97.54", "3.37", "0.0202", 6),
  @("516860", "博时中证金融科技主题ETF", "0.27", "98.47", "3.41", "0.0092", 6)
)

# Fund codes and the 基金规模/股票总仓位/仓位占比/持有市值 figures are stored
# as text in the source data (e.g. "008188" keeps its leading zero), so mark
# columns B:G as Text before writing them.
$newSheet.Range("B1:G16").NumberFormat = "@"

# Header row (B1:H1), bold / centered / boxed like every other sheet.
for ($col = 0; $col -lt $headers.Length; $col++) {
    $newSheet.Cells.Item(1, $col + 2).Value = $headers[$col]
}
$summary.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# Data rows (row 2..16).
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $row = $i + 2

    $newSheet.Cells.Item($row, 2).Value = $r[0]
    $newSheet.Cells.Item($row, 3).Value = $r[1]
    $newSheet.Cells.Item($row, 4).Value = $r[2]
    $newSheet.Cells.Item($row, 5).Value = $r[3]
    $newSheet.Cells.Item($row, 6).Value = $r[4]
    $newSheet.Cells.Item($row, 7).Value = $r[5]
    $newSheet.Cells.Item($row, 8).Value = $r[6]

    $newSheet.Cells.Item($row, 1).Value = $i
}

# Column A (the 0-based row index) uses the same boxed style as column A on
# every other sheet; copy it over from the summary sheet.
$summary.Range("A3").Copy()
$newSheet.Range("A2:A16").PasteSpecial(-4122)

$newSheet.Range("A1").Select()
